$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(2).ColumnWidth = 9.63
Write-Host "done"
